$p = $ppt.ActivePresentation

# --- Slide 15: "Possible Greedy Choices for Knapsack" ---
# Greedy choice #3 label: "profit-to-value" was wrong -> "value-to-weight"
$s15 = $p.Slides.Item(15)
$shp15 = $s15.Shapes.Item(3)
$shp15.TextFrame.TextRange.Runs(1).Text = "Greedy choice #3:  highest value-to-weight ratio"

# --- Slide 23: "Dynamic Programming" ---
# Corrected the (reversed) statement about when DP helps vs. greedy
$s23 = $p.Slides.Item(23)
$shp23 = $s23.Shapes.Item(3)
$para23 = $shp23.TextFrame.TextRange.Paragraphs(9)
$para23.Runs(1).Text = "DP is good when sub-problems overlap, when they’re not independent"

# --- Slide 7: "Example #1: Knapsack Problems" ---
# Drop the stray trailing period after "and a weight wi"
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(2)
$para7 = $shp7.TextFrame.TextRange.Paragraphs(2)
$para7.Runs(8).Text = " "
